# Rename the embedded logo pictures in the headers/footers of the first
# section, swapping image1.png <-> image2.png (Pearson footer logo) and
# image1.jpg <-> image2.jpg (BTEC header logo), matching the change made
# in Word's UI (e.g. Selection Pane "Rename" / re-saving the picture).
#
# wdHeaderFooterIndex: 1 = wdHeaderFooterPrimary, 2 = wdHeaderFooterFirstPage

$d = $word.ActiveDocument
$sec = $d.Sections.First

# --- Headers: BTEC logo, "image2.jpg" -> "image1.jpg"
for ($i = 1; $i -le 2; $i++) {
    $hdr = $sec.Headers.Item($i)
    if ($hdr.Exists) {
        $shapes = $hdr.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shapes.Item($j).Name = "image1.jpg"
        }
    }
}

# --- Footers: Pearson logo, "image1.png" -> "image2.png"
for ($i = 1; $i -le 2; $i++) {
    $ftr = $sec.Footers.Item($i)
    if ($ftr.Exists) {
        $shapes = $ftr.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shapes.Item($j).Name = "image2.png"
        }
    }
}
